$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Preserve A2's existing cell style (quotePrefix xf) before we overwrite its value ---
$ws.Range("A2").Copy()
$ws.Range("K1").PasteSpecial(-4122)   # xlPasteFormats -> stash format in helper cell K1

# --- Row 2 new content ---
$ws.Range("A2").Value = "Có công mài sắt có ngày nên ….."
$ws.Range("B2").Value = "kim"
$ws.Range("C2").Value = "Kim"
$ws.Range("D2").Value = 1
$ws.Range("E2").Value = "Người"
$ws.Range("F2").Value = 0
$ws.Range("G2").Value = "Dao"
$ws.Range("H2").Value = 0
$ws.Range("I2").Value = "Kéo"
$ws.Range("J2").Value = 0

# --- Restore A2's style (quotePrefix) that got wiped by the .Value assignment ---
$ws.Range("K1").Copy()
$ws.Range("A2").PasteSpecial(-4122)
$ws.Range("K1").Clear()

# --- Row 3 new content ---
$ws.Range("A3").Value = "Uống nước nhớ …."
$ws.Range("B3").Value = "Nguồn"
$ws.Range("C3").Value = "Nhà sản xuất"
$ws.Range("D3").Value = 0
$ws.Range("E3").Value = "Nguồn"
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = "Quả"
$ws.Range("H3").Value = 0
$ws.Range("I3").Value = "Cây"
$ws.Range("J3").Value = 0

# --- Row 4: clear out the old sample-question content; keep the styled-but-empty cells ---
$ws.Range("A4").ClearContents()
$ws.Range("B4").ClearContents()
$ws.Range("C4").ClearContents()
$ws.Range("D4").ClearContents()
$ws.Range("E4").ClearContents()
$ws.Range("F4").ClearContents()
$ws.Range("G4").ClearContents()
$ws.Range("H4").ClearContents()
$ws.Range("I4").ClearContents()
$ws.Range("J4").ClearContents()

# Row 4 shrinks back toward the default auto-fit height now that it's empty
$ws.Rows.Item(4).RowHeight = 17.4

# --- Selection moves to I4 (as last active cell) ---
$ws.Range("I4").Select()

$wb.Save()
